$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the B/C columns (thisletter / corrAns) for several rows ---
# Row 18: thisletter P -> X
$ws.Range("B18").Value = "X"

# Row 19: thisletter O -> X, corrAns 0 -> 1
$ws.Range("B19").Value = "X"
$ws.Range("C19").Value = 1

# Row 20: thisletter X -> P
$ws.Range("B20").Value = "P"

# Row 21: thisletter X -> O, corrAns 1 -> 0
$ws.Range("B21").Value = "O"
$ws.Range("C21").Value = 0

# Row 56: thisletter U -> Q
$ws.Range("B56").Value = "Q"

# Row 57: thisletter U -> A, corrAns 1 -> 0
$ws.Range("B57").Value = "A"
$ws.Range("C57").Value = 0

# Row 60: thisletter Q -> U
$ws.Range("B60").Value = "U"

# Row 61: thisletter A -> U, corrAns 0 -> 1
$ws.Range("B61").Value = "U"
$ws.Range("C61").Value = 1

# --- Add the new (blank, but formatted) D/E cells on rows 17-19, matching
#     the formatting already used by the A:C columns on those rows ---
$ws.Range("C17").Copy()
$ws.Range("D17:E17").PasteSpecial(-4122)

$ws.Range("C18").Copy()
$ws.Range("D18:E18").PasteSpecial(-4122)

$ws.Range("C19").Copy()
$ws.Range("D19:E19").PasteSpecial(-4122)

$ws.Range("C60").Copy()
$ws.Range("D60:E60").PasteSpecial(-4122)

$ws.Range("C61").Copy()
$ws.Range("D61:E61").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Move the active selection cursor, matching the saved cursor position ---
$ws.Range("F4").Select()
